# Update the "LamaMocogno" daily-cases sheet:
#   - a missing day (2021-02-08, serial 44235) is inserted at row 93,
#     pushing every following row down by one;
#   - the 7-day rolling-sum columns (C/D) for the now-complete week
#     ending 2021-02-27 (shifted to row 112) get their real values;
#   - two new trailing days (2021-03-01 / 2021-03-02) are appended as
#     rows 114/115.
# Net effect matches the diff: dimension grows from A1:D113 to A1:D115,
# and everything from row 93 on shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the two brand-new rows at the bottom (114/115, which don't
# exist yet) and the newly inserted row (93) inherit the same look as the
# rest of column A (style s="2": centred, bordered, date-formatted) before
# any value is poked into them - a bare ".Value =" on a previously empty
# cell does not inherit formatting.
$ws.Cells.Item(92, 1).Copy()
$ws.Cells.Item(114, 1).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(115, 1).PasteSpecial(-4122)   # xlPasteFormats

# Shift rows 93..113 down to 94..114 (walk bottom-up so we never clobber
# a source row before it has been read).
for ($r = 113; $r -ge 93; $r--) {
    $dst = $r + 1

    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2

    $cVal = $ws.Cells.Item($r, 3).Value2
    if ($cVal -ne $null -and $cVal -ne "") {
        $ws.Cells.Item($dst, 3).Value = $cVal
    }

    $dVal = $ws.Cells.Item($r, 4).Value2
    if ($dVal -ne $null -and $dVal -ne "") {
        $ws.Cells.Item($dst, 4).Value = $dVal
    }
}

# Row 93: the newly-inserted day (2021-02-08), zero new cases, style
# copied from the row above so column A keeps its date format/border.
$ws.Cells.Item(92, 1).Copy()
$ws.Cells.Item(93, 1).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(93, 1).Value = 44235
$ws.Cells.Item(93, 2).Value = 0
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 0

# Row 112 (2021-02-27, shifted down from the old row 111) now completes a
# full trailing 7-day window, so its rolling-sum columns get real numbers.
$ws.Cells.Item(112, 3).Value = 14
$ws.Cells.Item(112, 4).Value = 524.5410266017235

# New trailing row 115: 2021-03-02, zero new cases, C/D left blank.
$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 0
